$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Role in Project" row (row 3) with the new note text.
$ws.Range("C3").Value = "Xai bang Assigment role la Type"
$ws.Range("D3").Value = "0: ProjectOwer + PM, 1 : pM, 2: dev, 3: test, 4 QA, 5 cus, 6: project Owner`nProject Owner chi co' quyen read only va change PM,`nUser tao project se~ mang role la 0;"

# Wrap the long note text and make the row tall enough to show it.
$ws.Range("D3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 45

# Add the new "risk table" row.
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "risk table"
$ws.Range("C16").Value = "sudung table risk"

# Add the new "issue table" row.
$ws.Range("A17").Value = 12
$ws.Range("B17").Value = "issue table"
$ws.Range("C17").Value = "su dung table issue"
$ws.Range("D17").Value = "workUnitID will be projectID"

# Leave the selection on the cell that was edited, like Excel would.
$ws.Range("D3").Select() | Out-Null
